$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '51.604.58'
Set-TextValue 'E2' '  +1.17%  '
Set-TextValue 'D3' '3.022.14'
Set-TextValue 'E3' '  +2.40%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  +0.10%  '
Set-TextValue 'D5' '378.84'
Set-TextValue 'E5' '  -0.20%  '
Set-TextValue 'D6' '102.97'
Set-TextValue 'E6' '  +1.84%  '
Set-TextValue 'E7' '  +1.09%  '
Set-TextValue 'E8' '  +0.00%  '
Set-TextValue 'D9' '0.593'
Set-TextValue 'E9' '  +2.04%  '
Set-TextValue 'D10' '36.58'
Set-TextValue 'E12' '  +1.13%  '
Set-TextValue 'D13' '3.499.67'
Set-TextValue 'E13' '  +2.87%  '
Set-TextValue 'D14' '18.50'
Set-TextValue 'E14' '  +1.10%  '
Set-TextValue 'E15' '  +0.01%  '
Set-TextValue 'D16' '3.028.42'
Set-TextValue 'E16' '  +2.99%  '
Set-TextValue 'E17' '  -2.10%  '
Set-TextValue 'D18' '10.60'
Set-TextValue 'E18' '  -11.40%  '
Set-TextValue 'D19' '51.628.04'
Set-TextValue 'E19' '  +1.27%  '
Set-TextValue 'D20' '3.04'
Set-TextValue 'E20' '  -0.12%  '
Set-TextValue 'D21' '12.45'
Set-TextValue 'E21' '  +0.22%  '
Set-TextValue 'D22' '0.0₃0961'
Set-TextValue 'E22' '  +1.14%  '
Set-TextValue 'E23' '  +0.60%  '
Set-TextValue 'D24' '268.57'
Set-TextValue 'E24' '  +0.65%  '
Set-TextValue 'D25' '3.15'
Set-TextValue 'E25' '  -3.04%  '
Set-TextValue 'D26' '8.24'
Set-TextValue 'E26' '  +0.67%  '
Set-TextValue 'D27' '7.69'
Set-TextValue 'E27' '  +8.16%  '
Set-TextValue 'E28' '  +5.27%  '
Set-TextValue 'D29' '26.30'
Set-TextValue 'E29' '  +2.54%  '
Set-TextValue 'E30' '  +0.05%  '
Set-TextValue 'D31' '0.109'
Set-TextValue 'E31' '  +0.61%  '
Set-TextValue 'E32' '  +1.76%  '
Set-TextValue 'B33' 'InjectiveProtocol'
Set-TextValue 'C33' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D33' '34.14'
Set-TextValue 'E33' '  +1.84%  '
Set-TextValue 'B34' 'VeChain'
Set-TextValue 'C34' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D34' '0.0456'
Set-TextValue 'E34' '  +5.13%  '
Set-TextValue 'E35' '  +0.43%  '
Set-TextValue 'E37' '  -0.04%  '
Set-TextValue 'D38' '3.30'
Set-TextValue 'E38' '  +6.14%  '
Set-TextValue 'E39' '  +10.73%  '
Set-TextValue 'E40' '  +3.14%  '
Set-TextValue 'D41' '2.59'
Set-TextValue 'E41' '  +3.75%  '
Set-TextValue 'E42' '  +2.90%  '
Set-TextValue 'B43' 'Stellar'
Set-TextValue 'C43' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D43' '0.116'
Set-TextValue 'E43' '  -0.53%  '
Set-TextValue 'B44' 'Monero'
Set-TextValue 'C44' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D44' '126.75'
Set-TextValue 'E44' '  +5.98%  '
Set-TextValue 'D45' '3.74'
Set-TextValue 'E45' '  +7.64%  '
Set-TextValue 'E46' '  +0.99%  '
Set-TextValue 'E47' '  +3.45%  '
Set-TextValue 'D48' '2.38'
Set-TextValue 'E48' '  +2.44%  '
Set-TextValue 'D49' '2.031.37'
Set-TextValue 'E49' '  +0.96%  '
Set-TextValue 'D50' '3.323.62'
Set-TextValue 'E50' '  +2.66%  '
Set-TextValue 'D51' '0.0321'
Set-TextValue 'E51' '  +2.09%  '
